$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in hours invested and comments for the two new "curso de verano" entries
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = "● Curso de verano del moodle, capítulo 1"

$ws.Range("E16").Value = 1
$ws.Range("F16").Value = "● Curso de verano del moodle, capítulo 2"

# Underline formatting applied to F17 (left as-is, no content)
$ws.Range("F17").Font.Underline = $true

# Extend the day/date table down through row 22
$ws.Range("C18").Value = "Martes"
$ws.Range("D18").Value = 45363

$ws.Range("C19").Value = "Miércoles"
$ws.Range("D19").Value = 45364

$ws.Range("C20").Value = "Jueves"
$ws.Range("D20").Value = 45365

$ws.Range("C21").Value = "Viernes"
$ws.Range("D21").Value = 45366

$ws.Range("C22").Value = "Sábado"
$ws.Range("D22").Value = 45367

$ws.Range("E25").Select()
